# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on the Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6567.1665
$ws.Range("I76").Value = 7475.75
$ws.Range("J76").Value = 4750
$ws.Range("K76").Value = 7475.75
$ws.Range("L76").Value = 4750
$ws.Range("M76").Value = -7160.75
$ws.Range("N76").Value = -5380
$ws.Range("H79").Value = 6567.1665
$ws.Range("I79").Value = 7475.75
$ws.Range("J79").Value = 4750
$ws.Range("K79").Value = 7475.75
$ws.Range("L79").Value = 4750
$ws.Range("M79").Value = -6383.75
$ws.Range("N79").Value = -6934
$ws.Range("H80").Value = 5592.3184
$ws.Range("I80").Value = 756.1111
$ws.Range("J80").Value = 8940.462
$ws.Range("K80").Value = 2268.3333
$ws.Range("L80").Value = 26821.386
$ws.Range("M80").Value = -1270.3333
$ws.Range("N80").Value = -28817.386
$ws.Range("H83").Value = 5592.3184
$ws.Range("I83").Value = 756.1111
$ws.Range("J83").Value = 8940.462
$ws.Range("K83").Value = 6804.9999
$ws.Range("L83").Value = 80464.158
$ws.Range("M83").Value = -1812.9999
$ws.Range("N83").Value = -90448.158
$ws.Range("H86").Value = 3323.5833
$ws.Range("I86").Value = 2857.5715
$ws.Range("J86").Value = 3976
$ws.Range("K86").Value = 2857.5715
$ws.Range("L86").Value = 3976
$ws.Range("M86").Value = -1734.5715
$ws.Range("N86").Value = -6222
$ws.Range("H89").Value = 3323.5833
$ws.Range("I89").Value = 2857.5715
$ws.Range("J89").Value = 3976
$ws.Range("K89").Value = 14287.8575
$ws.Range("L89").Value = 19880
$ws.Range("M89").Value = -8671.8575
$ws.Range("N89").Value = -31112
$ws.Range("H92").Value = 1439
$ws.Range("I92").Value = 1392.9
$ws.Range("J92").Value = 1900
$ws.Range("K92").Value = 1392.9
$ws.Range("L92").Value = 1900
$ws.Range("M92").Value = -144.9000000000001
$ws.Range("N92").Value = -4396
$ws.Range("H138").Value = 2554.8572
$ws.Range("I138").Value = 1456.1765
$ws.Range("J138").Value = 3592.5
$ws.Range("K138").Value = 4368.529500000001
$ws.Range("L138").Value = 10777.5
$ws.Range("M138").Value = 771.4704999999994
$ws.Range("N138").Value = -21057.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 162.5
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 183.33333
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 183.33333
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = -407.33333
$ws.Range("H61").Value = 2902.0667
$ws.Range("I61").Value = 4043.6667
$ws.Range("J61").Value = 2141
$ws.Range("K61").Value = 4043.6667
$ws.Range("L61").Value = 2141
$ws.Range("M61").Value = -3831.6667
$ws.Range("N61").Value = -2565
$ws.Range("H132").Value = 6292.968
$ws.Range("I132").Value = 8375.058999999999
$ws.Range("K132").Value = 25125.177
$ws.Range("M132").Value = -22595.177
$ws.Range("H133").Value = 141148.67
$ws.Range("J133").Value = 141148.67
$ws.Range("L133").Value = 141148.67
$ws.Range("N133").Value = -146208.67
$ws.Range("H136").Value = 2902.0667
$ws.Range("I136").Value = 4043.6667
$ws.Range("J136").Value = 2141
$ws.Range("K136").Value = 12131.0001
$ws.Range("L136").Value = 6423
$ws.Range("M136").Value = -9581.000100000001
$ws.Range("N136").Value = -11523

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H4").Value = 162.5
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 183.33333
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 183.33333
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = -413.33333
$ws.Range("H134").Value = 2309.7778
$ws.Range("I134").Value = 1819
$ws.Range("J134").Value = 3291.3333
$ws.Range("K134").Value = 5457
$ws.Range("L134").Value = 9873.999899999999
$ws.Range("M134").Value = -2922
$ws.Range("N134").Value = -14943.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 28000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 28000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 28000
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -28348
$ws.Range("H22").Value = 359.48
$ws.Range("I22").Value = 212.11111
$ws.Range("J22").Value = 738.4286
$ws.Range("K22").Value = 212.11111
$ws.Range("L22").Value = 738.4286
$ws.Range("M22").Value = 137.88889
$ws.Range("N22").Value = -1438.4286
$ws.Range("H25").Value = 6190
$ws.Range("I25").Value = 4180
$ws.Range("J25").Value = 8200
$ws.Range("K25").Value = 4180
$ws.Range("L25").Value = 8200
$ws.Range("M25").Value = -4006
$ws.Range("N25").Value = -8548
$ws.Range("H132").Value = 502001.66
$ws.Range("I132").Value = 541681.9
$ws.Range("K132").Value = 1625045.7
$ws.Range("M132").Value = -1622515.7
$ws.Range("H134").Value = 1689.6061
$ws.Range("I134").Value = 1236.6522
$ws.Range("J134").Value = 2731.4
$ws.Range("K134").Value = 3709.9566
$ws.Range("L134").Value = 8194.200000000001
$ws.Range("M134").Value = -1174.9566
$ws.Range("N134").Value = -13264.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 56500
$ws.Range("J88").Value = 56500
$ws.Range("L88").Value = 169500
$ws.Range("N88").Value = -170356
$ws.Range("H91").Value = 56500
$ws.Range("J91").Value = 56500
$ws.Range("L91").Value = 169500
$ws.Range("N91").Value = -172464
$ws.Range("H137").Value = 2626.2104
$ws.Range("I137").Value = 1239.9
$ws.Range("J137").Value = 4166.5557
$ws.Range("K137").Value = 3719.7
$ws.Range("L137").Value = 12499.6671
$ws.Range("M137").Value = 1380.3
$ws.Range("N137").Value = -22699.6671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3493.2083
$ws.Range("I132").Value = 3308.3333
$ws.Range("J132").Value = 3801.3333
$ws.Range("K132").Value = 9924.999899999999
$ws.Range("L132").Value = 11403.9999
$ws.Range("M132").Value = -7394.999899999999
$ws.Range("N132").Value = -16463.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3616.9092
$ws.Range("I7").Value = 3385.75
$ws.Range("J7").Value = 4233.3335
$ws.Range("K7").Value = 3385.75
$ws.Range("L7").Value = 4233.3335
$ws.Range("M7").Value = -3273.75
$ws.Range("N7").Value = -4457.3335
$ws.Range("H126").Value = 3616.9092
$ws.Range("I126").Value = 3385.75
$ws.Range("J126").Value = 4233.3335
$ws.Range("K126").Value = 10157.25
$ws.Range("L126").Value = 12700.0005
$ws.Range("M126").Value = -7687.25
$ws.Range("N126").Value = -17640.0005
$ws.Range("H127").Value = 50715
$ws.Range("J127").Value = 50715
$ws.Range("L127").Value = 50715
$ws.Range("N127").Value = -60635
$ws.Range("H136").Value = 36073932
$ws.Range("I136").Value = 55557180
$ws.Range("K136").Value = 166671540
$ws.Range("M136").Value = -166668990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 60007
$ws.Range("J46").Value = 60007
$ws.Range("L46").Value = 60007
$ws.Range("N46").Value = -60469
$ws.Range("H119").Value = 275349
$ws.Range("J119").Value = 275349
$ws.Range("L119").Value = 275349
$ws.Range("N119").Value = -285025
$ws.Range("H132").Value = 1868.289
$ws.Range("I132").Value = 1196.1428
$ws.Range("J132").Value = 2975.353
$ws.Range("K132").Value = 3588.4284
$ws.Range("L132").Value = 8926.059000000001
$ws.Range("M132").Value = -1058.4284
$ws.Range("N132").Value = -13986.059
$ws.Range("H134").Value = 60007
$ws.Range("J134").Value = 60007
$ws.Range("L134").Value = 180021
$ws.Range("N134").Value = -185091
